$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (hunk 1)
$ws.Range("H12").Value = 286
$ws.Range("I12").Value = 66
$ws.Range("K12").Value = 66
$ws.Range("M12").Value = 104

# Row 92 (hunk 2)
$ws.Range("H92").Value = 1939.174
$ws.Range("I92").Value = 1780.45
$ws.Range("K92").Value = 1780.45
$ws.Range("M92").Value = -532.45

# Row 107 (hunk 3)
$ws.Range("H107").Value = 1727.4
$ws.Range("I107").Value = 1243.6
$ws.Range("K107").Value = 1243.6
$ws.Range("M107").Value = 676.4000000000001

# Row 112 (hunk 4)
$ws.Range("H112").Value = 2977.658
$ws.Range("J112").Value = 2921.3108
$ws.Range("L112").Value = 8763.932400000002
$ws.Range("N112").Value = -10979.9324

# Row 137 (hunk 5)
$ws.Range("H137").Value = 68969200
$ws.Range("I137").Value = 45457610
$ws.Range("K137").Value = 136372830
$ws.Range("M137").Value = -136370280

# Row 138 (hunk 6)
$ws.Range("H138").Value = 5736.579
$ws.Range("I138").Value = 2654
$ws.Range("J138").Value = 6314.5625
$ws.Range("K138").Value = 7962
$ws.Range("L138").Value = 18943.6875
$ws.Range("M138").Value = -2822
$ws.Range("N138").Value = -29223.6875

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 7)
$ws.Range("H32").Value = 14502336
$ws.Range("I32").Value = 18523652
$ws.Range("J32").Value = 25601
$ws.Range("K32").Value = 18523652
$ws.Range("L32").Value = 25601
$ws.Range("M32").Value = -18523365
$ws.Range("N32").Value = -26175

# Row 122 (hunk 8)
$ws.Range("H122").Value = 2570.8125
$ws.Range("I122").Value = 1706.8096
$ws.Range("J122").Value = 4220.273
$ws.Range("K122").Value = 5120.4288
$ws.Range("L122").Value = 12660.819
$ws.Range("M122").Value = -2670.4288
$ws.Range("N122").Value = -17560.819

# Row 132 (hunk 9)
$ws.Range("H132").Value = 23817770
$ws.Range("I132").Value = 9639.483
$ws.Range("J132").Value = 90913410
$ws.Range("K132").Value = 28918.449
$ws.Range("L132").Value = 272740230
$ws.Range("M132").Value = -26388.449
$ws.Range("N132").Value = -272745290

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (hunk 10)
$ws.Range("H86").Value = 16405.27
$ws.Range("I86").Value = 10124
$ws.Range("J86").Value = 30538.125
$ws.Range("K86").Value = 10124
$ws.Range("L86").Value = 30538.125
$ws.Range("M86").Value = -9001
$ws.Range("N86").Value = -32784.125

# Row 89 (hunk 11)
$ws.Range("H89").Value = 16405.27
$ws.Range("I89").Value = 10124
$ws.Range("J89").Value = 30538.125
$ws.Range("K89").Value = 50620
$ws.Range("L89").Value = 152690.625
$ws.Range("M89").Value = -45004
$ws.Range("N89").Value = -163922.625

# Row 105 (hunk 12)
$ws.Range("H105").Value = 12125.333
$ws.Range("I105").Value = 17529.166
$ws.Range("J105").Value = 1317.6666
$ws.Range("K105").Value = 17529.166
$ws.Range("L105").Value = 1317.6666
$ws.Range("M105").Value = -15782.166
$ws.Range("N105").Value = -4811.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 19 (hunk 13)
$ws.Range("H19").Value = 600.5454999999999
$ws.Range("I19").Value = 383.33334
$ws.Range("J19").Value = 861.2
$ws.Range("K19").Value = 383.33334
$ws.Range("L19").Value = 861.2
$ws.Range("M19").Value = -213.33334
$ws.Range("N19").Value = -1201.2

# Row 22 (hunk 14)
$ws.Range("H22").Value = 7032.25
$ws.Range("I22").Value = 14494.143
$ws.Range("J22").Value = 1228.5555
$ws.Range("K22").Value = 14494.143
$ws.Range("L22").Value = 1228.5555
$ws.Range("M22").Value = -14144.143
$ws.Range("N22").Value = -1928.5555

# Row 24 (hunk 15)
$ws.Range("H24").Value = 600.5454999999999
$ws.Range("I24").Value = 383.33334
$ws.Range("J24").Value = 861.2
$ws.Range("K24").Value = 383.33334
$ws.Range("L24").Value = 861.2
$ws.Range("M24").Value = -213.33334
$ws.Range("N24").Value = -1201.2

# Row 31 (hunk 16)
$ws.Range("H31").Value = 27031802
$ws.Range("I31").Value = 3850.6296
$ws.Range("K31").Value = 3850.6296
$ws.Range("M31").Value = -3555.6296

# Row 34 (hunk 17)
$ws.Range("H34").Value = 27031802
$ws.Range("I34").Value = 3850.6296
$ws.Range("K34").Value = 3850.6296
$ws.Range("M34").Value = -3648.6296

# Row 104 (hunk 18)
$ws.Range("H104").Value = 49964.332
$ws.Range("J104").Value = 49964.332
$ws.Range("L104").Value = 49964.332
$ws.Range("N104").Value = -55206.332

# Row 107 (hunk 19)
$ws.Range("H107").Value = 2703.7778
$ws.Range("I107").Value = 2337.6667
$ws.Range("J107").Value = 2886.8333
$ws.Range("K107").Value = 2337.6667
$ws.Range("L107").Value = 2886.8333
$ws.Range("M107").Value = -417.6667000000002
$ws.Range("N107").Value = -6726.8333

# Row 132 (hunk 20)
$ws.Range("H132").Value = 69256.734
$ws.Range("I132").Value = 85248.21000000001
$ws.Range("J132").Value = 5290.8335
$ws.Range("K132").Value = 255744.63
$ws.Range("L132").Value = 15872.5005
$ws.Range("M132").Value = -253214.63
$ws.Range("N132").Value = -20932.5005

# Row 141 (hunk 21)
$ws.Range("H141").Value = 302309.25
$ws.Range("J141").Value = 311904.47
$ws.Range("L141").Value = 311904.47
$ws.Range("N141").Value = -322264.47

$ws = $wb.Worksheets.Item("CUL")
# Row 92 (hunk 22)
$ws.Range("H92").Value = 194.25
$ws.Range("I92").Value = 194.25
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 582.75
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 665.25
$ws.Range("N92").ClearContents()

# Row 98 (hunk 23)
$ws.Range("H98").Value = 3057.5715
$ws.Range("I98").Value = 1663
$ws.Range("J98").Value = 3437.9092
$ws.Range("K98").Value = 4989
$ws.Range("L98").Value = 10313.7276
$ws.Range("M98").Value = -3491
$ws.Range("N98").Value = -13309.7276

# Row 107 (hunk 24)
$ws.Range("H107").Value = 598.44446
$ws.Range("J107").Value = 441.66666
$ws.Range("L107").Value = 1324.99998
$ws.Range("N107").Value = -5164.999980000001

# Row 113 (hunk 25)
$ws.Range("H113").Value = 3108.25
$ws.Range("I113").Value = 3699.5
$ws.Range("J113").Value = 3042.5557
$ws.Range("K113").Value = 11098.5
$ws.Range("L113").Value = 9127.667099999999
$ws.Range("M113").Value = -8928.5
$ws.Range("N113").Value = -13467.6671

# Row 128 (hunk 26)
$ws.Range("H128").Value = 116130
$ws.Range("I128").Value = 116130
$ws.Range("K128").Value = 348390
$ws.Range("M128").Value = -343410

# Row 131 (hunk 27)
$ws.Range("H131").Value = 44934.617
$ws.Range("J131").Value = 7365
$ws.Range("L131").Value = 22095
$ws.Range("N131").Value = -32175

# Row 136 (hunk 28)
$ws.Range("H136").Value = 5860.5
$ws.Range("I136").Value = 1721
$ws.Range("K136").Value = 5163
$ws.Range("M136").Value = -63

$ws = $wb.Worksheets.Item("GSM")
# Row 33 (hunk 29)
$ws.Range("H33").Value = 307298.8
$ws.Range("I33").Value = 5999
$ws.Range("J33").Value = 340776.56
$ws.Range("K33").Value = 5999
$ws.Range("L33").Value = 340776.56
$ws.Range("M33").Value = -5747
$ws.Range("N33").Value = -341280.56

# Row 70 (hunk 30)
$ws.Range("H70").Value = 60081.445
$ws.Range("I70").Value = 75665.53999999999
$ws.Range("J70").Value = 5537.125
$ws.Range("K70").Value = 75665.53999999999
$ws.Range("L70").Value = 5537.125
$ws.Range("M70").Value = -75395.53999999999
$ws.Range("N70").Value = -6077.125

# Row 73 (hunk 31)
$ws.Range("H73").Value = 60081.445
$ws.Range("I73").Value = 75665.53999999999
$ws.Range("J73").Value = 5537.125
$ws.Range("K73").Value = 75665.53999999999
$ws.Range("L73").Value = 5537.125
$ws.Range("M73").Value = -74729.53999999999
$ws.Range("N73").Value = -7409.125

# Row 104 (hunk 32)
$ws.Range("H104").Value = 15671
$ws.Range("J104").Value = 15671
$ws.Range("L104").Value = 15671
$ws.Range("N104").Value = -22659

# Row 113 (hunk 33)
$ws.Range("H113").Value = 5389.9287
$ws.Range("I113").Value = 5137.125
$ws.Range("K113").Value = 5137.125
$ws.Range("M113").Value = -2967.125

# Row 122 (hunk 34)
$ws.Range("H122").Value = 2213.2354
$ws.Range("I122").Value = 1893.8334
$ws.Range("J122").Value = 2979.8
$ws.Range("K122").Value = 5681.5002
$ws.Range("L122").Value = 8939.400000000001
$ws.Range("M122").Value = -3231.5002
$ws.Range("N122").Value = -13839.4

# Row 132 (hunk 35)
$ws.Range("H132").Value = 1978.1132
$ws.Range("I132").Value = 1875.6666
$ws.Range("J132").Value = 2369.2727
$ws.Range("K132").Value = 5626.9998
$ws.Range("L132").Value = 7107.8181
$ws.Range("M132").Value = -3096.9998
$ws.Range("N132").Value = -12167.8181

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (hunk 36)
$ws.Range("H40").Value = 5421.923
$ws.Range("I40").Value = 5317.364
$ws.Range("K40").Value = 5317.364
$ws.Range("M40").Value = -5181.364

# Row 46 (hunk 37)
$ws.Range("H46").Value = 1312.4
$ws.Range("I46").Value = 888.17645
$ws.Range("K46").Value = 888.17645
$ws.Range("M46").Value = -700.17645

# Row 122 (hunk 38)
$ws.Range("H122").Value = 5196.5
$ws.Range("I122").Value = 4593.4
$ws.Range("J122").Value = 5799.6
$ws.Range("K122").Value = 13780.2
$ws.Range("L122").Value = 17398.8
$ws.Range("M122").Value = -11330.2
$ws.Range("N122").Value = -22298.8

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (hunk 39)
$ws.Range("H54").Value = 12300
$ws.Range("J54").Value = 11000
$ws.Range("L54").Value = 11000
$ws.Range("N54").Value = -12040

# Row 107 (hunk 40)
$ws.Range("H107").Value = 642.5454999999999
$ws.Range("I107").Value = 585.3333
$ws.Range("K107").Value = 1755.9999
$ws.Range("M107").Value = 164.0001

# Row 113 (hunk 41)
$ws.Range("H113").Value = 953.5714
$ws.Range("I113").Value = 685
$ws.Range("J113").Value = 1061
$ws.Range("K113").Value = 2055
$ws.Range("L113").Value = 3183
$ws.Range("M113").Value = 115
$ws.Range("N113").Value = -7523
